$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.505.46'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.644.51'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '531.77'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.59'
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.569'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.68'
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.337'
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.115.91'
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '59.524.56'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.84'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.712.10'
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '343.93'
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.44'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.62'
$ws.Range('E20').Value = '  +2.70%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.39'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.18'
$ws.Range('E23').Value = '  +3.98%  '
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.769.78'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.20'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0801'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.38'
$ws.Range('E31').Value = '  -3.80%  '
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.05'
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '150.09'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.19'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.18'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.860'
$ws.Range('E37').Value = '  -4.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.855'
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.46'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.47'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.63'
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0980'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.601'
$ws.Range('E44').Value = '  -2.45%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '270.80'
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.39'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.72'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0535'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.039.99'
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.78'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  +0.99%  '
